$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 180; $r++) {
    $ws.Cells.Item($r, 5).Value = "V"
}
